$wb = $excel.ActiveWorkbook

# --- Text update: "Ready for handoff" -> "In Translation" everywhere it appears ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width update (report regenerated -> narrower Status-like columns) ---
# Target stored width ~13.41 chars; ColumnWidth snaps to the nearest achievable
# increment in this engine, so use the closest settable value.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
